$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 67
$ws1.Range("F3").Value = 798
$ws1.Range("F6").Value = 101
$ws1.Range("F7").Value = 334
$ws1.Range("F8").Value = 4079
$ws1.Range("F10").Value = 4797
$ws1.Range("F11").Value = 534
$ws1.Range("F12").Value = 1203

# Sheet "全部类型" (fourth sheet) - update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 67
$ws4.Range("F3").Value = 798
$ws4.Range("F6").Value = 101
$ws4.Range("F8").Value = 334
$ws4.Range("F9").Value = 4079
$ws4.Range("F11").Value = 4797
$ws4.Range("F12").Value = 534
$ws4.Range("F13").Value = 1203
